$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete old row 15 (Fonte Carregador Automotiva Jfa Storm 200a Monovolt Sci 220v) - removed in the edit
$ws.Rows(15).Delete()

# Step 2: insert two new columns at C:D for "modelo" and "politica"
$ws.Range("C1:D1").EntireColumn.Insert()

# Step 3: headers
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Step 5: fix B12 value (60.9 -> 60.8)
$ws.Range("B12").Value = 60.8

# Step 4: per-row data (modelo=C, politica=D, full=E, tipo=F, link=G); also fix tipo casing
# Row 2: Fonte Automotiva Jfa Storm 200a Bob Carregador Aut
$ws.Range("C2").Value = "FONTE 200 BOB"
$ws.Range("D2").Value = "Acima"
$ws.Range("F2").Value = "classico"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:579560342#searchVariation=MLB24834408&position=2&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 3: Fonte Carregador Jfa 70a Bivolt Com Medidor Cca
$ws.Range("C3").Value = "FONTE 70A"
$ws.Range("D3").Value = "Igual"
$ws.Range("F3").Value = "premium"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:579560342#searchVariation=MLB21455208&position=3&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 4: Controle Longa Distância Jfa Acqua 1200 Resistente
$ws.Range("C4").Value = "Sem Modelo"
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:579560342#searchVariation=MLB27687422&position=6&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 5: Fonte Carregador Bob Storm 200a Slim Bivolt Automá
$ws.Range("C5").Value = "FONTE 200 BOB"
$ws.Range("D5").Value = "Igual"
$ws.Range("F5").Value = "classico"
$ws.Range("G5").Value = "https://www.mercadolivre.com.br/fonte-carregador-bob-storm-200a-slim-bivolt-automatico-jfa-cor-preto/p/MLB27156459?pdp_filters=seller_id:579560342#searchVariation=MLB27156459&position=9&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 6: Controle Remoto Universal Longa Distância Jfa K120
$ws.Range("C6").Value = "Sem Modelo"
$ws.Range("D6").Value = ""
$ws.Range("F6").Value = "classico"
$ws.Range("G6").Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-vermelho/p/MLB34210379?pdp_filters=seller_id:579560342#searchVariation=MLB34210379&position=7&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 7: Controle Longa Distancia Jfa K1200 Alcance De 1200
$ws.Range("C7").Value = "Sem Modelo"
$ws.Range("D7").Value = ""
$ws.Range("F7").Value = "classico"
$ws.Range("G7").Value = "https://www.mercadolivre.com.br/controle-longa-distancia-jfa-k1200-alcance-de-1200-metros/p/MLB33922926?pdp_filters=seller_id:579560342#searchVariation=MLB33922926&position=4&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 8: Controle Longa Distancia Jfa K1200 Alcance De 1200
$ws.Range("C8").Value = "Sem Modelo"
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = "classico"
$ws.Range("G8").Value = "https://www.mercadolivre.com.br/controle-longa-distancia-jfa-k1200-alcance-de-1200-metros/p/MLB34245679?pdp_filters=seller_id:579560342#searchVariation=MLB34245679&position=5&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 9: Fonte Carregador Jfa 70 Amperes Slim Bivolt C/ Vol
$ws.Range("C9").Value = "FONTE 70A"
$ws.Range("D9").Value = "Acima"
$ws.Range("F9").Value = "classico"
$ws.Range("G9").Value = "https://produto.mercadolivre.com.br/MLB-2808437099-fonte-carregador-jfa-70-amperes-slim-bivolt-c-voltimetro-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 10: Fonte Carregador Jfa Bob Storm 200a Bivolt
$ws.Range("C10").Value = "FONTE 200 BOB"
$ws.Range("D10").Value = "Igual"
$ws.Range("F10").Value = "premium"
$ws.Range("G10").Value = "https://produto.mercadolivre.com.br/MLB-2753133396-fonte-carregador-jfa-bob-storm-200a-bivolt-_JM#position%3D11%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 11: Controle Longa Distancia Jfa K1200 Alcance 1200 Mt
$ws.Range("C11").Value = "Sem Modelo"
$ws.Range("D11").Value = ""
$ws.Range("F11").Value = "classico"
$ws.Range("G11").Value = "https://produto.mercadolivre.com.br/MLB-2715487599-controle-longa-distancia-jfa-k1200-alcance-1200-mt-_JM#position%3D12%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 12: Controle Longa Distancia Jfa K600 600 Metros Compl
$ws.Range("C12").Value = "Sem Modelo"
$ws.Range("D12").Value = ""
$ws.Range("F12").Value = "classico"
$ws.Range("G12").Value = "https://produto.mercadolivre.com.br/MLB-2069718298-controle-longa-distancia-jfa-k600-600-metros-completo-_JM#position%3D13%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 13: Controle Longa Distância Jfa Aqua 1200metros Branc
$ws.Range("C13").Value = "Sem Modelo"
$ws.Range("D13").Value = ""
$ws.Range("F13").Value = "premium"
$ws.Range("G13").Value = "https://produto.mercadolivre.com.br/MLB-2753158506-controle-longa-distncia-jfa-aqua-1200metros-branco-completo-_JM#position%3D14%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 14: Controle Longa Distancia Jfa Redline Wr P/aparelho
$ws.Range("C14").Value = "Sem Modelo"
$ws.Range("D14").Value = ""
$ws.Range("F14").Value = "premium"
$ws.Range("G14").Value = "https://produto.mercadolivre.com.br/MLB-2069705706-controle-longa-distancia-jfa-redline-wr-paparelho-original-_JM#position%3D15%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 15: Fonte Automotiva Jfa 200a Storm Voltímetro Digital
$ws.Range("C15").Value = "FONTE 200 MONO"
$ws.Range("D15").Value = "Igual"
$ws.Range("F15").Value = "classico"
$ws.Range("G15").Value = "https://produto.mercadolivre.com.br/MLB-3250326563-fonte-automotiva-jfa-200a-storm-voltimetro-digital-mono-220v-_JM#position%3D16%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 16: Filtro Anti-ruido Jfa Com Blindagem Eletromagnétic
$ws.Range("C16").Value = "Sem Modelo"
$ws.Range("D16").Value = ""
$ws.Range("F16").Value = "classico"
$ws.Range("G16").Value = "https://produto.mercadolivre.com.br/MLB-2139862487-filtro-anti-ruido-jfa-com-blindagem-eletromagnetica-2020k-_JM#position%3D17%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 17: Carregador De Bateria Carro Jfa 60a Redline Sci Po
$ws.Range("C17").Value = "Modelo identificado mas fora do range de preco"
$ws.Range("D17").Value = ""
$ws.Range("F17").Value = "classico"
$ws.Range("G17").Value = "https://produto.mercadolivre.com.br/MLB-2751443168-carregador-de-bateria-carro-jfa-60a-redline-sci-portatil-_JM#position%3D18%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 18: Carregador Jfa Storm 200a 12v 14,4v Smart Cca Sci 
$ws.Range("C18").Value = "FONTE 200 MONO"
$ws.Range("D18").Value = "Igual"
$ws.Range("F18").Value = "classico"
$ws.Range("G18").Value = "https://produto.mercadolivre.com.br/MLB-3246945241-carregador-jfa-storm-200a-12v-144v-smart-cca-sci-220v-mono-_JM#position%3D19%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 19: Controle Longa Distância Jfa Aqua 1200 Metros Pret
$ws.Range("C19").Value = "Sem Modelo"
$ws.Range("D19").Value = ""
$ws.Range("F19").Value = "premium"
$ws.Range("G19").Value = "https://produto.mercadolivre.com.br/MLB-2753141762-controle-longa-distncia-jfa-aqua-1200-metros-preto-completo-_JM#position%3D20%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 20: Fonte Carregador Jfa 200a Storm Voltímetro Digital
$ws.Range("C20").Value = "FONTE 200 MONO"
$ws.Range("D20").Value = "Igual"
$ws.Range("F20").Value = "classico"
$ws.Range("G20").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-storm-voltimetro-digital-mono-220v-cor-preto/p/MLB24006449?pdp_filters=seller_id:579560342#searchVariation=MLB24006449&position=1&search_layout=stack&type=product&tracking_id=fcb65d1f-f75f-446f-be9b-0fcf1857cef8"

# Row 21: Controle Longa Distancia Jfa Redline Wr P Aparelho
$ws.Range("C21").Value = "Sem Modelo"
$ws.Range("D21").Value = ""
$ws.Range("F21").Value = "classico"
$ws.Range("G21").Value = "https://produto.mercadolivre.com.br/MLB-2715542194-controle-longa-distancia-jfa-redline-wr-p-aparelho-original-_JM#position%3D21%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dfcb65d1f-f75f-446f-be9b-0fcf1857cef8"

